$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "31.418.20"
$ws.Range("E2").Value = "  +3.68%  "

# Row 3
$ws.Range("D3").Value = "2.007.79"
$ws.Range("E3").Value = "  +7.49%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7724"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +63.54%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "259.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.61%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3592"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +25.21%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.33"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +30.95%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07064"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.17%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8393"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +17.42%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08095"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.07%  "

# Row 13
$ws.Range("D13").Value = "2.007.48"
$ws.Range("E13").Value = "  +7.49%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "101.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.63%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.650"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +10.47%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.68%  "

# Row 17
$ws.Range("D17").Value = "31.420.31"
$ws.Range("E17").Value = "  +3.73%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +13.23%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.952"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +13.68%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007968"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.96%  "

# Row 21
$ws.Range("D21").Value = "2.269.82"
$ws.Range("E21").Value = "  +7.65%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.220"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +15.75%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.15%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1469"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +53.18%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.391"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +27.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.627"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.91%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.630"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.30%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.356"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.34%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.402"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.11%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05213"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.67%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.233"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.59%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7613"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.46%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.807"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.59%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02019"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.66%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.959"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.98%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.713"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.19%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "80.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.68%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.190"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.18%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4742"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.34%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8606"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.70%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.32%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.04%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4369"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.87%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.10%  "

# Row 28: only the volume % changes
$ws.Range("E28").Value = "  +7.92%  "

# Row 47: was EnergySwap, now Aptos (rows 47/48 swapped content)
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.672"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.85%  "

# Row 48: was Aptos, now EnergySwap
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.971"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.02%  "

# Row 51: was Maker, now Algorand
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1198"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +16.10%  "
